# Update the date label and all 26 multiplication problems in the table
# to the values for the "output generated at c8c62b6" commit.

$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $d.Content.Find.Execute($find, $false, $false, $false, $false, $false, `
                             $true, 1, $false, $replace, 2) | Out-Null
}

# Header date
Replace-Text "2025-09-06 Saturday" "2025-09-07 Sunday"

# Multiplication problems (row by row, left to right)
Replace-Text "49×20=" "87×18="
Replace-Text "80×23=" "70×60="
Replace-Text "77×40=" "78×73="
Replace-Text "93×14=" "27×51="
Replace-Text "49×56=" "38×59="

Replace-Text "87×91=" "57×77="
Replace-Text "61×25=" "71×16="
Replace-Text "20×90=" "59×60="
Replace-Text "48×41=" "95×53="
Replace-Text "62×36=" "45×81="

Replace-Text "41×37=" "31×43="
Replace-Text "71×62=" "39×78="
Replace-Text "37×17=" "30×55="
Replace-Text "67×97=" "58×47="
Replace-Text "29×22=" "95×71="

Replace-Text "83×41=" "57×98="
Replace-Text "52×70=" "13×44="
Replace-Text "13×92=" "16×88="
Replace-Text "21×63=" "33×14="
Replace-Text "58×78=" "87×41="

Replace-Text "71×33=" "25×91="
Replace-Text "74×51=" "46×46="
Replace-Text "64×84=" "98×32="
Replace-Text "56×22=" "27×69="
Replace-Text "77×60=" "32×40="
